$d = $word.ActiveDocument
$para = $d.Paragraphs(1)
$para.Range.Delete()
